$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 0.9083623842315833
$ws.Range("B2").Value = 103.3215974555178
$ws.Range("C2").Value = 0.003499286810108335
$ws.Range("D2").Value = 1285.529855385846
$ws.Range("E2").Value = 0.003216772622561443
$ws.Range("F2").Value = 1366.809244478769
$ws.Range("G2").Value = 1285.529905388346
$ws.Range("H2").Value = 2180.498742756034
$ws.Range("I2").Value = 0.5857522146917168
$ws.Range("J2").Value = 6.004058078522843
$ws.Range("K2").Value = 0.6776411608975116
$ws.Range("L2").Value = 1.171504429383434
$ws.Range("Y2").Value = 1265.382428459782
$ws.Range("Z2").Value = 267.9377834610459
$ws.Range("AA3").Value = 0.8710012221181069
$ws.Range("B3").Value = 103.3569733366944
$ws.Range("C3").Value = 0.002186786045805703
$ws.Range("D3").Value = 1285.483473993078
$ws.Range("E3").Value = 0.001945076639178862
$ws.Range("F3").Value = 1293.306052392802
$ws.Range("G3").Value = 1285.483523995578
$ws.Range("H3").Value = 2058.897395237981
$ws.Range("I3").Value = 0.5908860361644908
$ws.Range("J3").Value = 3.925972699705117
$ws.Range("K3").Value = 0.6511892080616624
$ws.Range("L3").Value = 1.181772072328982
$ws.Range("Y3").Value = 1265.30630909068
$ws.Range("Z3").Value = 230.3849896746154
$ws.Range("AA4").Value = 0.9079270719560177
$ws.Range("B4").Value = 103.3480953665314
$ws.Range("C4").Value = 0.00176474038680615
$ws.Range("D4").Value = 1285.494437632656
$ws.Range("E4").Value = 0.001549277845970626
$ws.Range("F4").Value = 1856.198326639334
$ws.Range("G4").Value = 1285.494487635156
$ws.Range("H4").Value = 2970.509853882965
$ws.Range("I4").Value = 0.5930636546058315
$ws.Range("J4").Value = 4.691197345167309
$ws.Range("K4").Value = 0.6550316697659154
$ws.Range("L4").Value = 1.186127309211663
$ws.Range("Y4").Value = 1265.246056570146
$ws.Range("Z4").Value = 354.1393799733564
$ws.Range("B5").Value = 103.0030537737121
$ws.Range("C5").Value = 0.03241960717311081
$ws.Range("D5").Value = 1286.359386528977
$ws.Range("E5").Value = 0.02944010052148074
$ws.Range("F5").Value = 46.4436303637359
$ws.Range("G5").Value = 1286.359386528977
$ws.Range("H5").Value = 46.28388218332019
$ws.Range("I5").Value = 0.4689848433978812
$ws.Range("J5").Value = 1.839396731452155
$ws.Range("K5").Value = 0.0000002031177993577771
$ws.Range("L5").Value = 0.9379696867957624
$ws.Range("AA6").Value = 0.5391480921773364
$ws.Range("B6").Value = 103.1550168005419
$ws.Range("C6").Value = 0.003336155775218801
$ws.Range("D6").Value = 1285.990416820285
$ws.Range("E6").Value = 0.002896860592757173
$ws.Range("F6").Value = 631.3070441520954
$ws.Range("G6").Value = 1285.990466822785
$ws.Range("H6").Value = 772.5818009097765
$ws.Range("I6").Value = 0.4822372590647193
$ws.Range("J6").Value = 2.542672903299207
$ws.Range("K6").Value = 0.4998668726168308
$ws.Range("L6").Value = 0.9644745181294385
$ws.Range("Y6").Value = 1265.695404865116
$ws.Range("Z6").Value = 79.48702084238347
$ws.Range("AA7").Value = 0.6841347470119548
$ws.Range("B7").Value = 103.1738792805686
$ws.Range("C7").Value = 0.004489142141953977
$ws.Range("D7").Value = 1285.955243558347
$ws.Range("E7").Value = 0.00384159600258908
$ws.Range("F7").Value = 415.3522796918753
$ws.Range("G7").Value = 1285.955293560847
$ws.Range("H7").Value = 534.6267116054648
$ws.Range("I7").Value = 0.5001689220274297
$ws.Range("J7").Value = 2.188066749627262
$ws.Range("K7").Value = 0.5360314218086061
$ws.Range("L7").Value = 1.000337844054859
$ws.Range("Y7").Value = 1265.650764676925
$ws.Range("Z7").Value = 56.9338881409605
$ws.Range("B8").Value = 103.0883319577022
$ws.Range("C8").Value = 0.007131543751507582
$ws.Range("D8").Value = 1286.164409157406
$ws.Range("E8").Value = 0.006425688134980971
$ws.Range("F8").Value = 235.2924360424752
$ws.Range("G8").Value = 1286.164409157406
$ws.Range("H8").Value = 270.8185567004009
$ws.Range("I8").Value = 0.4709534536529494
$ws.Range("J8").Value = 2.039371367791538
$ws.Range("K8").Value = 0.398927530251662
$ws.Range("L8").Value = 0.9419069073058988
$ws.Range("B9").Value = 103.3171074334812
$ws.Range("C9").Value = 0.006756076420931702
$ws.Range("D9").Value = 1285.61358880601
$ws.Range("E9").Value = 0.006324007328653321
$ws.Range("F9").Value = 323.8384992521369
$ws.Range("G9").Value = 1285.61358880601
$ws.Range("H9").Value = 490.5827186049489
$ws.Range("I9").Value = 0.5832211272049332
$ws.Range("J9").Value = 2.594740378412047
$ws.Range("K9").Value = 0.5619577483789917
$ws.Range("L9").Value = 1.166442254409866
$ws.Range("AA10").Value = 0.9069231919924228
$ws.Range("B10").Value = 103.3303943939397
$ws.Range("C10").Value = 0.001511149599471953
$ws.Range("D10").Value = 1285.56883093596
$ws.Range("E10").Value = 0.001255698845684445
$ws.Range("F10").Value = 1945.508074963094
$ws.Range("G10").Value = 1285.568880938461
$ws.Range("H10").Value = 3042.258120394317
$ws.Range("I10").Value = 0.5934214741934769
$ws.Range("J10").Value = 3.704510172281391
$ws.Range("K10").Value = 0.5962531282154446
$ws.Range("L10").Value = 1.186842948386954
$ws.Range("Y10").Value = 1265.321469083896
$ws.Range("Z10").Value = 359.0463679727656
$ws.Range("AA11").Value = 0.7738815825779851
$ws.Range("B11").Value = 103.2957472105443
$ws.Range("C11").Value = 0.001835957871401626
$ws.Range("D11").Value = 1285.608748390306
$ws.Range("E11").Value = 0.001589720073255292
$ws.Range("F11").Value = 1678.154312916877
$ws.Range("G11").Value = 1285.608798392806
$ws.Range("H11").Value = 2690.727819306308
$ws.Range("I11").Value = 0.6018353711801896
$ws.Range("J11").Value = 3.798007921122309
$ws.Range("K11").Value = 0.6235181888632202
$ws.Range("L11").Value = 1.203670742360379
$ws.Range("Y11").Value = 1265.337197657449
$ws.Range("Z11").Value = 303.0847703199496
$ws.Range("AA12").Value = 0.7957621290036879
$ws.Range("B12").Value = 103.2900608882007
$ws.Range("C12").Value = 0.003012197362429106
$ws.Range("D12").Value = 1285.637988404147
$ws.Range("E12").Value = 0.002670071333269119
$ws.Range("F12").Value = 620.692594105217
$ws.Range("G12").Value = 1285.638038406647
$ws.Range("H12").Value = 964.5465483537987
$ws.Range("I12").Value = 0.5920754740289971
$ws.Range("J12").Value = 2.421211449489129
$ws.Range("K12").Value = 0.5862116996535665
$ws.Range("L12").Value = 1.184150948057994
$ws.Range("Y12").Value = 1265.414010080515
$ws.Range("Z12").Value = 115.2448288101352
$ws.Range("AA13").Value = 0.6204214835512503
$ws.Range("B13").Value = 103.2834581875011
$ws.Range("C13").Value = 0.003528501496918605
$ws.Range("D13").Value = 1285.637300915657
$ws.Range("E13").Value = 0.003122445423862462
$ws.Range("F13").Value = 468.9697838099689
$ws.Range("G13").Value = 1285.637350918157
$ws.Range("H13").Value = 724.4512500196562
$ws.Range("I13").Value = 0.5991280012292467
$ws.Range("J13").Value = 1.889209250070597
$ws.Range("K13").Value = 0.5408858100621294
$ws.Range("L13").Value = 1.198256002458493
$ws.Range("Y13").Value = 1265.351937330548
$ws.Range("Z13").Value = 65.80178925310864
$ws.Range("AA14").Value = 0.7184263494221589
$ws.Range("B14").Value = 103.2095119051621
$ws.Range("C14").Value = 0.007089808265529515
$ws.Range("D14").Value = 1285.903623350436
$ws.Range("E14").Value = 0.006108657812127236
$ws.Range("F14").Value = 331.7713589142032
$ws.Range("G14").Value = 1285.903673352936
$ws.Range("H14").Value = 447.4899518263346
$ws.Range("I14").Value = 0.5429700163104516
$ws.Range("J14").Value = 2.547700141169557
$ws.Range("K14").Value = 0.4435646080716731
$ws.Range("L14").Value = 1.085940032620903
$ws.Range("Y14").Value = 1265.685710086654
$ws.Range("Z14").Value = 52.25416031540784
$ws.Range("AA15").Value = 0.5458776516798269
$ws.Range("B15").Value = 103.2520857124941
$ws.Range("C15").Value = 0.008550769370825509
$ws.Range("D15").Value = 1285.82039257487
$ws.Range("E15").Value = 0.007370482204753549
$ws.Range("F15").Value = 667.9908014468608
$ws.Range("G15").Value = 1285.82044257737
$ws.Range("H15").Value = 936.6403205264056
$ws.Range("I15").Value = 0.5681795407454743
$ws.Range("J15").Value = 6.244208513802181
$ws.Range("K15").Value = 0.4262206857897317
$ws.Range("L15").Value = 1.136359081490949
$ws.Range("Y15").Value = 1265.462056221352
$ws.Range("Z15").Value = 109.2723011597471
$ws.Range("AA16").Value = 1.045940250198919
$ws.Range("B16").Value = 103.2167621413046
$ws.Range("C16").Value = 0.01726630379687906
$ws.Range("D16").Value = 1285.894349070691
$ws.Range("E16").Value = 0.01439774970445403
$ws.Range("F16").Value = 406.3179622280051
$ws.Range("G16").Value = 1285.894399073191
$ws.Range("H16").Value = 568.0998052337924
$ws.Range("I16").Value = 0.5601547531387823
$ws.Range("J16").Value = 8.591087799189003
$ws.Range("K16").Value = 0.4569736898901452
$ws.Range("L16").Value = 1.120309506277565
$ws.Range("Y16").Value = 1265.823228252162
$ws.Range("Z16").Value = 131.2439190594631
$ws.Range("AA17").Value = 0.8501745208316925
$ws.Range("B17").Value = 103.2845203905702
$ws.Range("C17").Value = 0.006010913181031239
$ws.Range("D17").Value = 1285.692316838596
$ws.Range("E17").Value = 0.005444931735024882
$ws.Range("F17").Value = 477.3951566135385
$ws.Range("G17").Value = 1285.692366841096
$ws.Range("H17").Value = 795.8622986585715
$ws.Range("I17").Value = 0.6701336724756033
$ws.Range("J17").Value = 3.378069679130222
$ws.Range("K17").Value = 0.4477814916158721
$ws.Range("L17").Value = 1.340267344951207
$ws.Range("Y17").Value = 1265.513590940886
$ws.Range("Z17").Value = 126.2247705647952

$ws.Range("AA8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()

Write-Output "Applied diff changes"
